# Generate Report for Handoff
# Refresh the localization-status report for the newly generated source
# file: the old GUID-named file (070629b4-0e14-4a13-8590-da5bd76ca603)
# is replaced throughout by the newly generated one
# (dd1ade58-150e-4c56-95f1-83f5ac8e691b), and the handoff timestamps /
# xliff hash tokens are refreshed to match the new handoff run.

$wb = $excel.ActiveWorkbook

$newGuid = "dd1ade58-150e-4c56-95f1-83f5ac8e691b"
$newHash = "32c4d075845bc27d193646e8f56b189d8ec69212"

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33ed0d2d3cb0375936063f542307ac6d08189797/e2e"

# ---------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "$repoBlobBase/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$newGuid.md"
)

$wsOverview.Range("G2").Value = "2016-08-19 23:02:58"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "$repoBlobBase/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
)

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-19 23:02:54"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "$repoBlobBase/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
)

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-19 23:02:58"
